{"js": "// Move the \"LOB1021 -  F\u00edsica IV  (Requisito fraco)\" requirement line\n// (text + manual line break) from the start of the \"Requisitos\" bullet\n// paragraph to the end, after \"LOM3257 -  Mec\u00e2nica Cl\u00e1ssica\".\n\nconst LOB_TEXT = \"LOB1021 -  F\u00edsica IV  (Requisito fraco)\";\nconst ANCHOR_TEXT = \"LOM3253\";\n\n// Find the paragraph that contains the LOB1021 requirement line.\nconst searchResults = context.document.body.search(LOB_TEXT, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the LOB1021 requirement text.\");\n}\n\nconst hit = searchResults.items[0];\nconst para = hit.paragraphs.getFirst();\npara.load(\"text\");\nawait context.sync();\n\n// Locate the anchor (\"LOM3253 ...\") that immediately follows the LOB1021\n// run in the original paragraph, so we can capture the exact span\n// (text + trailing line break) that needs to move.\nconst anchorResults = para.search(ANCHOR_TEXT, { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error(\"Could not find the LOM3253 anchor text.\");\n}\n\nconst anchorStart = anchorResults.items[0].getRange(\"Start\");\nconst paraStart = para.getRange(\"Start\");\n\n// The span from the paragraph start up to (but excluding) the anchor is\n// exactly \"LOB1021 -  F\u00edsica IV  (Requisito fraco)\" followed by the\n// manual line break (\\u000b) that separates it from the next entry.\nconst movingSpan = paraStart.expandTo(anchorStart);\nmovingSpan.load(\"text\");\nawait context.sync();\n\nconst movingText = movingSpan.text;\n\n// Remove the span from its current (first) position...\nmovingSpan.delete();\nawait context.sync();\n\n// ...and re-insert the same text (including its trailing manual line\n// break) at the end of the paragraph, after \"LOM3257 ...\".\nconst paraEnd = para.getRange(\"End\");\nparaEnd.insertText(movingText, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Move the \"LOB1021 -  F\u00edsica IV  (Requisito fraco)\" requirement line\n# (text + manual line break) from the start of the \"Requisitos\" bullet\n# paragraph to the end, after \"LOM3257 -  Mec\u00e2nica Cl\u00e1ssica\".\n\n$d = $word.ActiveDocument\n\n$LOB_TEXT = \"LOB1021 -  F\u00edsica IV  (Requisito fraco)\"\n$ANCHOR_TEXT = \"LOM3253\"\n\n# Find the paragraph that contains the LOB1021 requirement line by\n# walking the (live) Paragraphs collection, so the Range we keep stays\n# attached to the paragraph and auto-adjusts as we edit it below.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $cand = $d.Paragraphs.Item($i)\n    if ($cand.Range.Text.Contains($LOB_TEXT)) {\n        $target = $cand\n        break\n    }\n}\nif ($null -eq $target) {\n    throw \"Could not find the LOB1021 requirement text.\"\n}\n$pRange = $target.Range\n\n# Locate the anchor (\"LOM3253 ...\") that immediately follows the LOB1021\n# run in the original paragraph, so we can capture the exact span\n# (text + trailing manual line break) that needs to move.\n$anchorRange = $pRange.Duplicate\n$anchorFound = $anchorRange.Find.Execute($ANCHOR_TEXT)\nif (-not $anchorFound) {\n    throw \"Could not find the LOM3253 anchor text.\"\n}\n$anchorStart = $anchorRange.Start\n\n# The span from the paragraph start up to (but excluding) the anchor is\n# exactly \"LOB1021 -  F\u00edsica IV  (Requisito fraco)\" followed by the\n# manual line break character that separates it from the next entry.\n$movingSpan = $d.Range($pRange.Start, $anchorStart)\n$movingText = $movingSpan.Text\n\n# Remove the span from its current (first) position. $pRange is the\n# paragraph's own (live) range, so its End automatically shrinks to\n# reflect the deletion.\n$movingSpan.Delete() | Out-Null\n\n# ...and re-insert the same text (including its trailing manual line\n# break) at the end of the paragraph, i.e. right before the paragraph\n# mark, so it lands after \"LOM3257 ...\".\n$pEnd = $pRange.End\n$insertPoint = $d.Range($pEnd - 1, $pEnd - 1)\n$insertPoint.InsertAfter($movingText) | Out-Null\n"}
